# API: update new customer
# Regenerate the "Password" column (column E) values for the customer rows
# and move the active-cell selection down one row, matching the upstream
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

$passwords = @{
    2  = "VJbbNX"
    3  = "GPhNHj"
    4  = "xvbOIa"
    5  = "ZbQUkX"
    6  = "leqNom"
    7  = "icYwig"
    9  = "saHVAu"
    10 = "HYfdiH"
    11 = "wshKan"
    12 = "BnsJhU"
    13 = "sBuGyF"
}

foreach ($row in $passwords.Keys) {
    $ws.Cells.Item($row, 5).Value = $passwords[$row]
}

$ws.Range("E16").Select()
